# Update results summary file
# - rename sheet from "CEC HW Main Meter Power" to "CEC HW Main Meter Flow"
# - change the dataset filename referenced in column A (rows 2-4) from the
#   "1a" csv to the "3a" csv
# - move the active selection to A5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "CEC HW Main Meter Flow"

# Update the "Dataset" column values (rows 2-4) to the new CSV file name
$ws.Range("A2").Value = "CEC_compiled_data_3a_updated.csv"
$ws.Range("A3").Value = "CEC_compiled_data_3a_updated.csv"
$ws.Range("A4").Value = "CEC_compiled_data_3a_updated.csv"

# Move the selection to A5 (matches saved cursor position in the file)
$ws.Range("A5").Select()
